# Tasks Details.xlsx - update per commit "All task till 27-Dec-2020 07:52 pm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G width (new column for "Payment Status") ---
$ws.Columns.Item(7).ColumnWidth = 17.85546875

# --- Header row ---
$ws.Range("G1").Value = "Payment Status"

# --- Row 2: Excel and DB Task ---
$ws.Range("E2").Value = "Delivered"

# --- Row 3: Java and DB Task ---
$ws.Range("E3").Value = "Delivered"

# --- Row 4: Java Thred Task (unchanged values, left as-is) ---

# --- Row 5: PostgreSQL (unchanged values, left as-is) ---

# --- Row 6: Database for Ticket System -- re-themed from "blue" to "green" block,
#            and work status corrected from "In Progress" to "Dilivered" ---
$ws.Range("A6:G6").Interior.Color = 5296274
$ws.Range("A6:G6").HorizontalAlignment = -4108
$ws.Range("A6:G6").VerticalAlignment = -4108
$ws.Range("C6:D6").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

$ws.Range("E6").Value = "Dilivered"
$ws.Range("G6").Value = "Done"

# --- Row 7: new task "Excel to Access Conversion" (blue themed, like old row 6) ---
$ws.Range("A7:G7").Interior.Color = 15773696
$ws.Range("A7:G7").HorizontalAlignment = -4108
$ws.Range("A7:G7").VerticalAlignment = -4108
$ws.Range("C7:D7").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

$ws.Range("A7").Value = "Excel to Access Conversion"
$ws.Range("B7").Value = 1
$ws.Range("C7").Formula = "=TODAY()"
$ws.Range("D7").Formula = "=TODAY()"
$ws.Range("E7").Value = "Dilivered"
$ws.Range("F7").Value = "1k"
$ws.Range("G7").Value = "Pending"

# --- Selection as last left by the editor ---
$ws.Range("C15").Select()
